# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2 through 23
$kValues = @{
    2  = 0
    3  = 4
    4  = 3
    5  = 5
    6  = 3
    7  = 4
    8  = 7
    9  = 5
    10 = 6
    11 = 3
    12 = 2
    13 = 1
    14 = 1
    15 = 1
    16 = 2
    17 = 1
    18 = 0
    19 = 1
    20 = 0
    21 = 0
    22 = 1
    23 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
